$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "68.200.18"
$ws.Cells.Item(2, 5).Value = "  +0.01%  "
$ws.Cells.Item(3, 4).Value = "3.663.40"
$ws.Cells.Item(3, 5).Value = "  -1.15%  "
$ws.Cells.Item(4, 5).Value = "  +0.16%  "
$c = $ws.Cells.Item(5, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "594.87"
$c.Style = $origStyle
$ws.Cells.Item(5, 5).Value = "  -0.44%  "
$c = $ws.Cells.Item(6, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "191.63"
$c.Style = $origStyle
$ws.Cells.Item(6, 5).Value = "  +4.84%  "
$c = $ws.Cells.Item(7, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.621"
$c.Style = $origStyle
$ws.Cells.Item(7, 5).Value = "  -1.09%  "
$ws.Cells.Item(8, 5).Value = "  +0.52%  "
$c = $ws.Cells.Item(9, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.697"
$c.Style = $origStyle
$ws.Cells.Item(9, 5).Value = "  -2.64%  "
$c = $ws.Cells.Item(10, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.152"
$c.Style = $origStyle
$ws.Cells.Item(10, 5).Value = "  -6.52%  "
$c = $ws.Cells.Item(11, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "56.92"
$c.Style = $origStyle
$ws.Cells.Item(11, 5).Value = "  +1.36%  "
$c = $ws.Cells.Item(12, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.0000271"
$c.Style = $origStyle
$ws.Cells.Item(12, 5).Value = "  -6.69%  "
$ws.Cells.Item(13, 5).Value = "  -1.32%  "
$ws.Cells.Item(14, 4).Value = "4.255.03"
$ws.Cells.Item(14, 5).Value = "  -0.80%  "
$ws.Cells.Item(15, 4).Value = "3.665.82"
$ws.Cells.Item(15, 5).Value = "  -0.91%  "
$ws.Cells.Item(16, 5).Value = "  +0.06%  "
$c = $ws.Cells.Item(17, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "18.86"
$c.Style = $origStyle
$ws.Cells.Item(17, 5).Value = "  -2.94%  "
$ws.Cells.Item(18, 2).Value = "WrappedBTC"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Cells.Item(18, 4).Value = "68.053.41"
$ws.Cells.Item(18, 5).Value = "  -0.05%  "
$ws.Cells.Item(19, 2).Value = "Uniswap"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$c = $ws.Cells.Item(19, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "12.60"
$c.Style = $origStyle
$ws.Cells.Item(19, 5).Value = "  -1.81%  "
$ws.Cells.Item(20, 2).Value = "Polygon"
$ws.Cells.Item(20, 3).Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$c = $ws.Cells.Item(20, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.10"
$c.Style = $origStyle
$ws.Cells.Item(20, 5).Value = "  -1.96%  "
$c = $ws.Cells.Item(21, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "405.20"
$c.Style = $origStyle
$ws.Cells.Item(21, 5).Value = "  -1.10%  "
$c = $ws.Cells.Item(22, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "4.43"
$c.Style = $origStyle
$ws.Cells.Item(22, 5).Value = "  -3.03%  "
$c = $ws.Cells.Item(23, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "88.39"
$c.Style = $origStyle
$ws.Cells.Item(23, 5).Value = "  -0.33%  "
$c = $ws.Cells.Item(24, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "2.94"
$c.Style = $origStyle
$ws.Cells.Item(24, 5).Value = "  -2.51%  "
$c = $ws.Cells.Item(25, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "10.96"
$c.Style = $origStyle
$ws.Cells.Item(25, 5).Value = "  -0.66%  "
$c = $ws.Cells.Item(26, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "12.52"
$c.Style = $origStyle
$ws.Cells.Item(26, 5).Value = "  -2.26%  "
$c = $ws.Cells.Item(27, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "6.08"
$c.Style = $origStyle
$ws.Cells.Item(27, 5).Value = "  +0.23%  "
$c = $ws.Cells.Item(28, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "3.67"
$c.Style = $origStyle
$ws.Cells.Item(28, 5).Value = "  -4.74%  "
$c = $ws.Cells.Item(29, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "9.31"
$c.Style = $origStyle
$ws.Cells.Item(29, 5).Value = "  -1.95%  "
$c = $ws.Cells.Item(30, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "31.87"
$c.Style = $origStyle
$ws.Cells.Item(30, 5).Value = "  -2.80%  "
$c = $ws.Cells.Item(31, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "7.15"
$c.Style = $origStyle
$ws.Cells.Item(31, 5).Value = "  -1.62%  "
$c = $ws.Cells.Item(32, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "68.14"
$c.Style = $origStyle
$ws.Cells.Item(32, 5).Value = "  +6.27%  "
$c = $ws.Cells.Item(33, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "12.26"
$c.Style = $origStyle
$ws.Cells.Item(33, 5).Value = "  -2.02%  "
$c = $ws.Cells.Item(34, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "44.16"
$c.Style = $origStyle
$ws.Cells.Item(34, 5).Value = "  +1.24%  "
$ws.Cells.Item(35, 5).Value = "  -0.88%  "
$c = $ws.Cells.Item(36, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "606.42"
$c.Style = $origStyle
$ws.Cells.Item(36, 5).Value = "  +2.22%  "
$ws.Cells.Item(37, 5).Value = "  -0.13%  "
$c = $ws.Cells.Item(38, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.391"
$c.Style = $origStyle
$ws.Cells.Item(38, 5).Value = "  -2.52%  "
$ws.Cells.Item(39, 5).Value = "  -0.09%  "
$ws.Cells.Item(40, 4).Value = "0.0₃0770"
$ws.Cells.Item(40, 5).Value = "  -12.75%  "
$ws.Cells.Item(41, 5).Value = "  -0.82%  "
$c = $ws.Cells.Item(42, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "2.89"
$c.Style = $origStyle
$ws.Cells.Item(42, 5).Value = "  -3.83%  "
$c = $ws.Cells.Item(43, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.0425"
$c.Style = $origStyle
$ws.Cells.Item(43, 5).Value = "  -2.40%  "
$c = $ws.Cells.Item(44, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "2.51"
$c.Style = $origStyle
$ws.Cells.Item(44, 5).Value = "  -9.00%  "
$c = $ws.Cells.Item(45, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "3.25"
$c.Style = $origStyle
$ws.Cells.Item(45, 5).Value = "  +3.87%  "
$c = $ws.Cells.Item(46, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.135"
$c.Style = $origStyle
$ws.Cells.Item(46, 5).Value = "  +0.98%  "
$ws.Cells.Item(47, 4).Value = "2.777.80"
$ws.Cells.Item(47, 5).Value = "  +0.23%  "
$c = $ws.Cells.Item(48, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "8.90"
$c.Style = $origStyle
$ws.Cells.Item(48, 5).Value = "  -3.35%  "
$c = $ws.Cells.Item(49, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "143.66"
$c.Style = $origStyle
$ws.Cells.Item(49, 5).Value = "  +1.81%  "
$c = $ws.Cells.Item(50, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "2.58"
$c.Style = $origStyle
$ws.Cells.Item(50, 5).Value = "  -5.08%  "
$c = $ws.Cells.Item(51, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "2.58"
$c.Style = $origStyle
$ws.Cells.Item(51, 5).Value = "  -11.07%  "
